# Insert a new weekly record as row 71 on "Hortaliza, Terminal La Palmera de
# La Serena - Cilantro" sheet, pushing the existing rows 71..107 down to
# 72..108 (dimension grows from A1:R107 to A1:R108).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data (rows 71-107) down by one row, inserting a blank row.
$ws.Rows.Item(71).Insert()

# Populate the newly inserted row 71 with the new weekly record.
$ws.Range("A71").Value = 8
$ws.Range("B71").Value = "Terminal La Palmera de La Serena"
$ws.Range("C71").Value = "Coquimbo"
$ws.Range("D71").Value = 44553
$ws.Range("E71").Value = 4
$ws.Range("F71").Value = 100112040
$ws.Range("G71").Value = "Cilantro"
$ws.Range("H71").Value = "Sin especificar"
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 3300
$ws.Range("K71").Value = 2000
$ws.Range("L71").Value = 2500
$ws.Range("M71").Value = 2250
$ws.Range("N71").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O71").Value = "Provincia del Elquí"
$ws.Range("P71").Value = 1500
$ws.Range("Q71").Value = 1.5
$ws.Range("R71").Value = "Hortaliza"
